$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 1: unmerge the old grouped header bands and give every column its own header text ---
$ws.Range("H1:L1").UnMerge()
$ws.Range("M1:P1").UnMerge()
$ws.Range("Q1:S1").UnMerge()

$ws.Range("A1").Value = "Player ID"
$ws.Range("B1").Value = "Player"
$ws.Range("C1").Value = "#"
$ws.Range("D1").Value = "Nation"
$ws.Range("E1").Value = "Pos"
$ws.Range("F1").Value = "Age"
$ws.Range("G1").Value = "90s"
$ws.Range("H1").Value = "Tkl"
$ws.Range("I1").Value = "TklW"
$ws.Range("J1").Value = "Def 3rd"
$ws.Range("K1").Value = "Mid 3rd"
$ws.Range("L1").Value = "Att 3rd"
$ws.Range("M1").Value = "Cha"
$ws.Range("N1").Value = "Att"
$ws.Range("O1").Value = "Tkl%"
$ws.Range("P1").Value = "Lost"
$ws.Range("Q1").Value = "Blocks"
$ws.Range("R1").Value = "Sh"
$ws.Range("S1").Value = "Pass"
$ws.Range("T1").Value = "Int"
$ws.Range("U1").Value = "Tkl+Int"
$ws.Range("V1").Value = "Clr"
$ws.Range("W1").Value = "Err"

# --- Row 2 becomes the hidden "real" sub-header row ---
$ws.Rows.Item(2).Hidden = $true

# --- A brand-new blank hidden spacer row at position 3 ---
$ws.Rows.Item(3).Hidden = $true

# --- A few rows were missing an explicit 0 in the Tkl% (O) column; make it explicit ---
$ws.Range("O5").Value = 0
$ws.Range("O8").Value = 0
$ws.Range("O10").Value = 0
$ws.Range("O19").Value = 0

# --- Totals row (20) gets hidden and its Tkl% carries full float precision ---
$ws.Rows.Item(20).Hidden = $true
$ws.Range("O20").Value = 39.299999999999997
